$wb = $excel.ActiveWorkbook

# Update the selection on the "Functional_diversity" sheet (L16)
$wsFunctional = $wb.Worksheets.Item("Functional_diversity")
$wsFunctional.Activate()
$wsFunctional.Range("L16").Select()

# Finally activate "Social_fabric" so it becomes the active/selected tab
$wsSocial = $wb.Worksheets.Item("Social_fabric")
$wsSocial.Activate()
